# "fix excel add template"
# Replace the old single-cell placeholder (A20: "Template Mẫu ") with the
# two-line report footer text, anchored back at the top of the sheet
# (A1/A2), and leave the selection where the author left off (H4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old placeholder cell so it doesn't linger down at row 20.
$ws.Range("A20").ClearContents()

# New template content.
$ws.Range("A1").Value = "KẾT THÚC BÁO CÁO"
$ws.Range("A2").Value = "CÓ VẤN ĐỀ GÌ HÃY BÁO LẠI QUA EMAIL :bachh1124@gmail.com"

# Match the saved selection/active cell from the authored workbook.
$ws.Range("H4").Select() | Out-Null
